$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

$data = New-Object 'object[,]' 50,4

$data[0,0] = 'Bitcoin'
$data[0,1] = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$data[0,2] = '25.808.44'
$data[0,3] = '  -5.17%  '
$data[1,0] = 'Ethereum'
$data[1,1] = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$data[1,2] = '1.813.11'
$data[1,3] = '  -4.37%  '
$data[2,0] = 'TetherUSD'
$data[2,1] = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$data[2,2] = '0.9994'
$data[2,3] = '  -0.27%  '
$data[3,0] = 'BNB'
$data[3,1] = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$data[3,2] = '276.64'
$data[3,3] = '  -9.66%  '
$data[4,0] = 'USDC'
$data[4,1] = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$data[4,2] = '0.9996'
$data[4,3] = '  -0.23%  '
$data[5,0] = 'XRP'
$data[5,1] = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$data[5,2] = '0.5104'
$data[5,3] = '  -5.20%  '
$data[6,0] = 'Cardano'
$data[6,1] = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$data[6,2] = '0.3521'
$data[6,3] = '  -7.14%  '
$data[7,0] = 'OKB'
$data[7,1] = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$data[7,2] = '44.56'
$data[7,3] = '  -2.47%  '
$data[8,0] = 'Dogecoin'
$data[8,1] = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$data[8,2] = '0.06668'
$data[8,3] = '  -8.37%  '
$data[9,0] = 'Solana'
$data[9,1] = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$data[9,2] = '20.11'
$data[9,3] = '  -8.58%  '
$data[10,0] = 'Polygon'
$data[10,1] = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$data[10,2] = '0.8347'
$data[10,3] = '  -6.98%  '
$data[11,0] = 'TRON'
$data[11,1] = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$data[11,2] = '0.07823'
$data[11,3] = '  -4.34%  '
$data[12,0] = 'WrappedEther'
$data[12,1] = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$data[12,2] = '1.796.07'
$data[12,3] = '  +0.37%  '
$data[13,0] = 'Polkadot'
$data[13,1] = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$data[13,2] = '5.077'
$data[13,3] = '  -4.99%  '
$data[14,0] = 'Litecoin'
$data[14,1] = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$data[14,2] = '87.98'
$data[14,3] = '  -7.29%  '
$data[15,0] = 'BinanceUSD'
$data[15,1] = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$data[15,2] = '0.9990'
$data[15,3] = '  -0.32%  '
$data[16,0] = 'Avalanche'
$data[16,1] = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$data[16,2] = '13.90'
$data[16,3] = '  -6.30%  '
$data[17,0] = 'ShibaInu'
$data[17,1] = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$data[17,2] = '0.000008020'
$data[17,3] = '  -7.20%  '
$data[18,0] = 'Dai'
$data[18,1] = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$data[18,2] = '0.9988'
$data[18,3] = '  -0.29%  '
$data[19,0] = 'WrappedBTC'
$data[19,1] = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$data[19,2] = '25.865.54'
$data[19,3] = '  -4.11%  '
$data[20,0] = 'Uniswap'
$data[20,1] = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$data[20,2] = '4.731'
$data[20,3] = '  -6.02%  '
$data[21,0] = 'Cosmos'
$data[21,1] = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$data[21,2] = '10.01'
$data[21,3] = '  -7.29%  '
$data[22,0] = 'Chainlink'
$data[22,1] = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$data[22,2] = '6.061'
$data[22,3] = '  -6.52%  '
$data[23,0] = 'Monero'
$data[23,1] = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$data[23,2] = '141.85'
$data[23,3] = '  -4.60%  '
$data[24,0] = 'LidoDAOToken'
$data[24,1] = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$data[24,2] = '2.201'
$data[24,3] = '  -3.87%  '
$data[25,0] = 'Toncoin'
$data[25,1] = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$data[25,2] = '1.654'
$data[25,3] = '  -5.90%  '
$data[26,0] = 'EthereumClassic'
$data[26,1] = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$data[26,2] = '17.05'
$data[26,3] = '  -7.08%  '
$data[27,0] = 'BitcoinCash'
$data[27,1] = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$data[27,2] = '108.85'
$data[27,3] = '  -6.30%  '
$data[28,0] = 'InternetComputer(DFINITY)'
$data[28,1] = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$data[28,2] = '4.339'
$data[28,3] = '  -9.92%  '
$data[29,0] = 'Filecoin'
$data[29,1] = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$data[29,2] = '4.225'
$data[29,3] = '  -9.04%  '
$data[30,0] = 'Stellar'
$data[30,1] = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$data[30,2] = '0.08788'
$data[30,3] = '  -4.08%  '
$data[31,0] = 'Hedera'
$data[31,1] = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$data[31,2] = '0.04886'
$data[31,3] = '  -3.09%  '
$data[32,0] = 'ImmutableX'
$data[32,1] = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$data[32,2] = '0.7330'
$data[32,3] = '  -10.59%  '
$data[33,0] = 'ARBITRUM'
$data[33,1] = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$data[33,2] = '1.140'
$data[33,3] = '  -6.44%  '
$data[34,0] = 'HuobiToken'
$data[34,1] = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$data[34,2] = '2.891'
$data[34,3] = '  -4.28%  '
$data[35,0] = 'Frax'
$data[35,1] = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$data[35,2] = '0.9986'
$data[35,3] = '  -0.17%  '
$data[36,0] = 'MXToken'
$data[36,1] = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$data[36,2] = '3.048'
$data[36,3] = '  -7.43%  '
$data[37,0] = 'TheSandbox'
$data[37,1] = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$data[37,2] = '0.5229'
$data[37,3] = '  -12.22%  '
$data[38,0] = 'VeChain'
$data[38,1] = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$data[38,2] = '0.01857'
$data[38,3] = '  -6.42%  '
$data[39,0] = 'RenderToken'
$data[39,1] = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$data[39,2] = '2.295'
$data[39,3] = '  -14.13%  '
$data[40,0] = 'TrustWalletToken'
$data[40,1] = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$data[40,2] = '0.9538'
$data[40,3] = '  -11.23%  '
$data[41,0] = 'Quant'
$data[41,1] = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$data[41,2] = '112.09'
$data[41,3] = '  -2.38%  '
$data[42,0] = 'FraxShare'
$data[42,1] = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$data[42,2] = '6.183'
$data[42,3] = '  -6.66%  '
$data[43,0] = 'Aptos'
$data[43,1] = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$data[43,2] = '8.125'
$data[43,3] = '  -11.99%  '
$data[44,0] = 'PaxDollar'
$data[44,1] = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$data[44,2] = '0.9990'
$data[44,3] = '  -0.23%  '
$data[45,0] = 'Decentraland'
$data[45,1] = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$data[45,2] = '0.4570'
$data[45,3] = '  -10.13%  '
$data[46,0] = 'Algorand'
$data[46,1] = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$data[46,2] = '0.1382'
$data[46,3] = '  -9.41%  '
$data[47,0] = 'EnergySwap'
$data[47,1] = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$data[47,2] = '9.319'
$data[47,3] = '  -8.30%  '
$data[48,0] = 'Elrond'
$data[48,1] = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$data[48,2] = '36.22'
$data[48,3] = '  -4.46%  '
$data[49,0] = 'NEARProtocol'
$data[49,1] = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$data[49,2] = '1.501'
$data[49,3] = '  -7.78%  '

$rng.Value = $data
